$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Narrow a set of columns from width 8 to width 7 (custom accuracy reformat) ---
# ColumnWidth is expressed in "characters"; Excel stores the column width in the
# sheet XML as characters + 5/6 (padding in MDW units), so subtract 5/6 to land
# exactly on the target stored width of 7.
$narrowCols = @(2,3,7,11,12,15,22,24,27,29,30,34)
foreach ($c in $narrowCols) {
    $ws.Columns($c).ColumnWidth = 7 - (5/6)
}

# --- Reduce row 5 values to 2 decimal places ("custom accuracy") ---
$row5Values = @(
    16.4,     # B5
    12,       # C5
    0.75,     # D5
    34.83,    # E5
    28.32,    # F5
    12.15,    # G5
    46.29,    # H5
    18.96,    # I5
    8.38,     # J5
    12.59,    # K5
    13.84,    # L5
    15.18,    # M5
    3.98,     # N5
    11.99,    # O5
    17.71,    # P5
    10.24,    # Q5
    0.44,     # R5
    0.36,     # S5
    181.33,   # T5
    34.61,    # U5
    11.64,    # V5
    23.59,    # W5
    12.49,    # X5
    1.62,     # Y5
    22.58,    # Z5
    10.11,    # AA5
    8.59,     # AB5
    10.53,    # AC5
    14.67,    # AD5
    0.12,     # AE5
    41.31,    # AF5
    6.37,     # AG5
    14.15     # AH5
)
for ($i = 0; $i -lt $row5Values.Length; $i++) {
    $ws.Cells.Item(5, $i + 2).Value = $row5Values[$i]
}

# --- Remove row 6 entirely (data trimmed to a single sample row) ---
$ws.Rows(6).Delete()
